# Applies the "design document" touch-up: Word re-ran its proofer
# (spelling/grammar squiggles + the auto _GoBack bookmark) and two
# literal " x" placeholders were typed at the end of two list items.
#
# Because proofErr/bookmark markers and run-splits aren't exposed as
# discrete Word object-model calls, each affected paragraph's Range is
# replaced wholesale via Range.InsertXML with the exact WordprocessingML
# we need (same pPr, split into the runs the diff shows).

$d = $word.ActiveDocument
$W = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Set-ParaXml($index, $innerXml) {
    $p = $d.Paragraphs.Item($index)
    $wrapped = "<w:p " + $W + ">" + $innerXml + "</w:p>"
    $p.Range.InsertXML($wrapped) | Out-Null
}

$listPPr = "<w:pPr><w:pStyle w:val='Listeafsnit'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr>"

# 1) "Bingeflix" title -> wrap in a spell-check span
$idx = 1
$xml = "<w:pPr><w:pStyle w:val='Titel'/></w:pPr>" `
    + "<w:proofErr w:type='spellStart'/>" `
    + "<w:r><w:t>Bingeflix</w:t></w:r>" `
    + "<w:proofErr w:type='spellEnd'/>"
Set-ParaXml $idx $xml

# 6) "Dato for første udgivelse" -> append literal " x" as its own run
$idx = 6
$xml = $listPPr `
    + "<w:r><w:t>Dato for f&#248;rste udgivelse</w:t></w:r>" `
    + "<w:r><w:t xml:space='preserve'> x</w:t></w:r>"
Set-ParaXml $idx $xml

# 7) "Dato for seneste udgivelse" -> append literal " x" as its own run
$idx = 7
$xml = $listPPr `
    + "<w:r><w:t>Dato for seneste udgivelse</w:t></w:r>" `
    + "<w:r><w:t xml:space='preserve'> x</w:t></w:r>"
Set-ParaXml $idx $xml

# 10) "Færdig eller on-going" -> split so "going" is spell-flagged
$idx = 10
$xml = $listPPr `
    + "<w:r><w:t>F&#230;rdig eller on-</w:t></w:r>" `
    + "<w:proofErr w:type='spellStart'/>" `
    + "<w:r><w:t>going</w:t></w:r>" `
    + "<w:proofErr w:type='spellEnd'/>"
Set-ParaXml $idx $xml

# 11) "Original-sprog" -> wrap in a spell-check span
$idx = 11
$xml = $listPPr `
    + "<w:proofErr w:type='spellStart'/>" `
    + "<w:r><w:t>Original-sprog</w:t></w:r>" `
    + "<w:proofErr w:type='spellEnd'/>"
Set-ParaXml $idx $xml

# 13) "Aldersgrænse" -> split around the relocated _GoBack bookmark
$idx = 13
$xml = $listPPr `
    + "<w:r><w:t>Aldersgr</w:t></w:r>" `
    + "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" `
    + "<w:bookmarkEnd w:id='0'/>" `
    + "<w:r><w:t>&#230;nse</w:t></w:r>"
Set-ParaXml $idx $xml

# 18) "IMDB / Rotten tomatoes rating" -> flag "tomatoes"
$idx = 18
$xml = $listPPr `
    + "<w:r><w:t xml:space='preserve'>IMDB / Rotten </w:t></w:r>" `
    + "<w:proofErr w:type='spellStart'/>" `
    + "<w:r><w:t>tomatoes</w:t></w:r>" `
    + "<w:proofErr w:type='spellEnd'/>" `
    + "<w:r><w:t xml:space='preserve'> rating</w:t></w:r>"
Set-ParaXml $idx $xml

# 21) "Når brugeren opretter sin profil vælger de " -> flag "profil" as a
# grammar issue; the paragraph's second run is untouched.
$idx = 21
$xml = "<w:r><w:t xml:space='preserve'>N&#229;r brugeren opretter sin </w:t></w:r>" `
    + "<w:proofErr w:type='gramStart'/>" `
    + "<w:r><w:t>profil</w:t></w:r>" `
    + "<w:proofErr w:type='gramEnd'/>" `
    + "<w:r><w:t xml:space='preserve'> v&#230;lger de </w:t></w:r>" `
    + "<w:r><w:t>en r&#230;kke pr&#230;ferencer:</w:t></w:r>"
Set-ParaXml $idx $xml

# 23) "Om de foretrækker et bestemt original-sprog" -> flag "original-sprog"
$idx = 23
$xml = $listPPr `
    + "<w:r><w:t xml:space='preserve'>Om de foretr&#230;kker et bestemt </w:t></w:r>" `
    + "<w:proofErr w:type='spellStart'/>" `
    + "<w:r><w:t>original-sprog</w:t></w:r>" `
    + "<w:proofErr w:type='spellEnd'/>"
Set-ParaXml $idx $xml

# 30) "Vi vil bruge SQLite da appen ..." -> flag "SQLite"; also the old
# trailing _GoBack bookmark here is superseded by the one added in
# paragraph 13 above (Word keeps only one _GoBack).
$idx = 30
$xml = "<w:r><w:t xml:space='preserve'>Vi vil bruge </w:t></w:r>" `
    + "<w:proofErr w:type='spellStart'/>" `
    + "<w:r><w:t>SQLite</w:t></w:r>" `
    + "<w:proofErr w:type='spellEnd'/>" `
    + "<w:r><w:t xml:space='preserve'> da appen </w:t></w:r>" `
    + "<w:r><w:t>ikke skal udgives rigtigt og der er derfor ingen grund til at ops&#230;tte en ekstern database server.</w:t></w:r>"
Set-ParaXml $idx $xml
